$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$days = @("Sunday","Monday","Tuesday","Wednesday","Thursday","Friday","Saturday")
$months = @("january ","february","march","april","may","june","july","august","september","october","november","december")

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $days[$i % 7]
}

for ($i = 0; $i -lt 12; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $months[$i]
}

for ($i = 0; $i -lt 12; $i++) {
    $ws.Cells.Item($i + 1, 3).Value = 2
}

$null = $ws.Range("C1:C12").Select()
